$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "부산대 오류수정": re-code the A column (계열/division) for every admissions
# row. The workbook previously stored each row's specific college name
# (인문대학, 공과대학, 약학대학, ...); the corrected sheet instead uses the
# broader 3-way exam-track classification used for 수능(가/나) grouping:
# 인문사회계 / 자연계 / 예술계.
$ws.Range("A2:A4").Value = "자연계"
$ws.Range("A5").Value = "인문사회계"
$ws.Range("A6").Value = "자연계"
$ws.Range("A7:A9").Value = "인문사회계"
$ws.Range("A10").Value = "자연계"
$ws.Range("A11:A12").Value = "인문사회계"
$ws.Range("A13").Value = "자연계"
$ws.Range("A14:A16").Value = "인문사회계"
$ws.Range("A17:A19").Value = "자연계"
$ws.Range("A20").Value = "인문사회계"
$ws.Range("A21").Value = "자연계"
$ws.Range("A22").Value = "인문사회계"
$ws.Range("A23").Value = "자연계"
$ws.Range("A24").Value = "예술계"
$ws.Range("A25:A26").Value = "인문사회계"
$ws.Range("A27:A28").Value = "자연계"
$ws.Range("A29").Value = "인문사회계"
$ws.Range("A30:A34").Value = "자연계"
$ws.Range("A35:A39").Value = "인문사회계"
$ws.Range("A40:A48").Value = "자연계"
$ws.Range("A49:A52").Value = "인문사회계"
$ws.Range("A53").Value = "자연계"
$ws.Range("A54:A58").Value = "인문사회계"
$ws.Range("A59:A60").Value = "자연계"
$ws.Range("A61:A63").Value = "인문사회계"
$ws.Range("A64:A65").Value = "자연계"
$ws.Range("A66:A67").Value = "인문사회계"
$ws.Range("A68:A72").Value = "자연계"
$ws.Range("A73").Value = "인문사회계"
$ws.Range("A74:A75").Value = "자연계"
$ws.Range("A76").Value = "인문사회계"
$ws.Range("A77").Value = "자연계"
$ws.Range("A78").Value = "인문사회계"
$ws.Range("A79").Value = "자연계"
$ws.Range("A80").Value = "인문사회계"
$ws.Range("A81:A82").Value = "자연계"
$ws.Range("A83:A84").Value = "인문사회계"
$ws.Range("A85:A87").Value = "자연계"
$ws.Range("A88").Value = "인문사회계"
$ws.Range("A89:A93").Value = "자연계"

# Leave the cursor/scroll position where the author left it after making
# these edits (selection on B85, default top-left cell).
$excel.Goto($ws.Range("B85"))
